$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "São Paulo"
$ws.Range("B11").Value = "26/08/2025 21:05"

# C11/D11 look numeric ("16"/"90") but must be stored as text, like the
# rest of the sheet. Force text number format, assign, then strip the
# format back off so no style index is left on the cell.
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "16"
$ws.Range("C11").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "90"
$ws.Range("D11").ClearFormats()

$ws.Range("E11").Value = "Nublado"
